$wb = $excel.ActiveWorkbook

# --- Update the Password values (shortened) on the approver sheets ---
$firstLevel = $wb.Worksheets.Item("FirstLevelApprover")
$firstLevel.Range("B2").Value = "Bingo@12345"
$firstLevel.Range("B3").Value = "Bingo@12345"
$firstLevel.Range("B4").Value = "Bingo@12345"
$firstLevel.Range("B5").Value = "Bingo@12345"

$approver = $wb.Worksheets.Item("Approver")
$approver.Range("B2").Value = "Bingo@12345"
$approver.Range("B3").Value = "Bingo@12345"
$approver.Range("B4").Value = "Bingo@12345"
$approver.Range("B5").Value = "Bingo@12345"

# --- Update the Approver sheet's remembered selection ---
$approver.Range("B6").Select()

# --- Make FirstLevelApprover the active/selected tab with its remembered selection ---
$firstLevel.Activate()
$firstLevel.Range("H19").Select()
